# Insert a new data row at row 8 (pushing existing rows 8-33 down to 9-34),
# then populate the new row with the new "Primera" / Región Metropolitana entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = 11
$ws.Range("B8").Value = "Vega Monumental Concepción"
$ws.Range("C8").Value = "Bíobío"
$ws.Range("D8").Value = 44659
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 100112037
$ws.Range("G8").Value = "Cebollín"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 8000
$ws.Range("L8").Value = 8500
$ws.Range("M8").Value = 8250
$ws.Range("N8").Value = "`$/paquete 36 unidades"
$ws.Range("O8").Value = "Región Metropolitana"
$ws.Range("P8").Value = 229
$ws.Range("Q8").Value = 36
$ws.Range("R8").Value = "Hortaliza"
